$wb = $excel.ActiveWorkbook

# --- Rename third sheet: "Reject app by Interpreter" -> "Reject-Accept app" ---
$wsReject = $wb.Worksheets.Item("Reject app by Interpreter")
$wsReject.Name = "Reject-Accept app"

# --- Sheet "New appointment": add column M (Testname / Ashwini.) ---
$wsNew = $wb.Worksheets.Item("New appointment")

# Copy formatting from the existing header cell (L1) onto the new header cell (M1)
$wsNew.Range("L1").Copy()
$wsNew.Range("M1").PasteSpecial(-4122) # xlPasteFormats

$wsNew.Cells.Item(1, 13).Value = "Testname"
$wsNew.Cells.Item(2, 13).Value = "Ashwini."

# Selection on "New appointment" moves to the whole first row, and it is no
# longer the tab-selected sheet.
$wsNew.Activate()
$wsNew.Range("A1:XFD1").Select()

# --- "Reject-Accept app" becomes the active / tab-selected sheet ---
$wsReject.Activate()
$wsReject.Range("D23").Select()

$wb.Save()
